$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new text value looks like a plain number need to be pre-formatted
# as Text so Excel keeps storing them as strings (matching the source data,
# which stores every Price/Volume column as text).
$textCells = @("D5", "D7", "D9", "D10", "D12", "D13", "D15", "D19", "D22", "D26", "D27", "D28", "D29", "D30", "D33", "D34", "D35", "D36", "D39", "D40", "D42", "D45", "D48", "D49", "D51")
foreach ($addr in $textCells) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range('D2').Value = '37.309.53'
$ws.Range('E2').Value = '  +2.03%  '
$ws.Range('D3').Value = '2.035.64'
$ws.Range('E3').Value = '  +3.84%  '
$ws.Range('E4').Value = '  +0.06%  '
$ws.Range('D5').Value = '248.32'
$ws.Range('E5').Value = '  +2.01%  '
$ws.Range('E6').Value = '  +2.31%  '
$ws.Range('D7').Value = '60.72'
$ws.Range('E7').Value = '  +0.21%  '
$ws.Range('E8').Value = '  +0.04%  '
$ws.Range('D9').Value = '0.399'
$ws.Range('E9').Value = '  +6.14%  '
$ws.Range('D10').Value = '0.0813'
$ws.Range('E10').Value = '  +3.07%  '
$ws.Range('E11').Value = '  +2.13%  '
$ws.Range('D12').Value = '15.38'
$ws.Range('E12').Value = '  +8.17%  '
$ws.Range('D13').Value = '0.865'
$ws.Range('E13').Value = '  +4.17%  '
$ws.Range('D14').Value = '2.335.44'
$ws.Range('E14').Value = '  +3.89%  '
$ws.Range('D15').Value = '22.44'
$ws.Range('E15').Value = '  +2.82%  '
$ws.Range('E16').Value = '  +5.26%  '
$ws.Range('D17').Value = '2.027.64'
$ws.Range('E17').Value = '  +3.50%  '
$ws.Range('D18').Value = '37.249.26'
$ws.Range('E18').Value = '  +2.12%  '
$ws.Range('D19').Value = '70.89'
$ws.Range('E19').Value = '  +1.82%  '
$ws.Range('E20').Value = '  +2.15%  '
$ws.Range('E21').Value = '  +3.88%  '
$ws.Range('D22').Value = '231.27'
$ws.Range('E22').Value = '  +0.72%  '
$ws.Range('E23').Value = '  +0.13%  '
$ws.Range('E24').Value = '  +3.47%  '
$ws.Range('E25').Value = '  +0.67%  '
$ws.Range('D26').Value = '9.50'
$ws.Range('E26').Value = '  +3.33%  '
$ws.Range('D27').Value = '163.77'
$ws.Range('E27').Value = '  +1.49%  '
$ws.Range('D28').Value = '0.136'
$ws.Range('E28').Value = '  -2.85%  '
$ws.Range('D29').Value = '19.89'
$ws.Range('E29').Value = '  +2.86%  '
$ws.Range('D30').Value = '1.40'
$ws.Range('E30').Value = '  +5.83%  '
$ws.Range('E31').Value = '  +2.85%  '
$ws.Range('E32').Value = '  +2.21%  '
$ws.Range('D33').Value = '0.0671'
$ws.Range('E33').Value = '  +9.33%  '
$ws.Range('B34').Value = 'LidoDAOToken'
$ws.Range('C34').Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range('D34').Value = '2.53'
$ws.Range('E34').Value = '  +11.43%  '
$ws.Range('B35').Value = 'InternetComputer(DFINITY)'
$ws.Range('C35').Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range('D35').Value = '4.54'
$ws.Range('E35').Value = '  +1.84%  '
$ws.Range('D36').Value = '3.57'
$ws.Range('E36').Value = '  +2.64%  '
$ws.Range('E37').Value = '  -0.01%  '
$ws.Range('E38').Value = '  +2.04%  '
$ws.Range('D39').Value = '5.44'
$ws.Range('E39').Value = '  -0.61%  '
$ws.Range('D40').Value = '0.0983'
$ws.Range('E41').Value = '  +1.52%  '
$ws.Range('D42').Value = '17.26'
$ws.Range('E42').Value = '  +9.05%  '
$ws.Range('E43').Value = '  +2.73%  '
$ws.Range('E44').Value = '  +2.91%  '
$ws.Range('D45').Value = '93.01'
$ws.Range('E45').Value = '  +4.73%  '
$ws.Range('E46').Value = '  +3.92%  '
$ws.Range('D47').Value = '1.390.02'
$ws.Range('E47').Value = '  +2.01%  '
$ws.Range('D48').Value = '7.54'
$ws.Range('E48').Value = '  +6.07%  '
$ws.Range('D49').Value = '2.18'
$ws.Range('E49').Value = '  +19.79%  '
$ws.Range('E50').Value = '  +1.19%  '
$ws.Range('D51').Value = '46.68'
$ws.Range('E51').Value = '  +1.94%  '
